$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H43").Value = 1681.25
$ws.Range("I43").Value = 1642.8572
$ws.Range("K43").Value = 1642.8572
$ws.Range("M43").Value = -1573.8572

$ws.Range("H70").Value = 126835.875
$ws.Range("J70").Value = 126835.875
$ws.Range("L70").Value = 380507.625
$ws.Range("N70").Value = -381047.625

$ws.Range("H73").Value = 126835.875
$ws.Range("J73").Value = 126835.875
$ws.Range("L73").Value = 380507.625
$ws.Range("N73").Value = -382379.625

$ws.Range("H92").Value = 200
$ws.Range("I92").Value = 200
$ws.Range("K92").Value = 200
$ws.Range("M92").Value = 1048

$ws.Range("H132").Value = 1338.4048
$ws.Range("I132").Value = 1440.4595
$ws.Range("K132").Value = 4321.3785
$ws.Range("M132").Value = -1791.3785

$ws = $wb.Worksheets("ARM")
$ws.Range("H32").Value = 5143.776
$ws.Range("I32").Value = 5143.776
$ws.Range("K32").Value = 5143.776
$ws.Range("M32").Value = -4856.776

$ws.Range("H63").Value = 5415.25
$ws.Range("I63").Value = 2805.6667
$ws.Range("K63").Value = 2805.6667
$ws.Range("M63").Value = -2119.6667

$ws.Range("H66").Value = 5415.25
$ws.Range("I66").Value = 2805.6667
$ws.Range("K66").Value = 14028.3335
$ws.Range("M66").Value = -10596.3335

$ws.Range("H74").Value = 1351.6316
$ws.Range("I74").Value = 1351.6316
$ws.Range("K74").Value = 1351.6316
$ws.Range("M74").Value = -477.6315999999999

$ws.Range("H77").Value = 1351.6316
$ws.Range("I77").Value = 1351.6316
$ws.Range("K77").Value = 6758.157999999999
$ws.Range("M77").Value = -2390.157999999999

$ws.Range("H88").Value = 3324.5
$ws.Range("I88").Value = 4646.8335
$ws.Range("J88").Value = 2002.1666
$ws.Range("K88").Value = 4646.8335
$ws.Range("L88").Value = 2002.1666
$ws.Range("M88").Value = -4240.8335
$ws.Range("N88").Value = -2814.1666

$ws.Range("H91").Value = 3324.5
$ws.Range("I91").Value = 4646.8335
$ws.Range("J91").Value = 2002.1666
$ws.Range("K91").Value = 4646.8335
$ws.Range("L91").Value = 2002.1666
$ws.Range("M91").Value = -3242.8335
$ws.Range("N91").Value = -4810.1666

$ws = $wb.Worksheets("BSM")
$ws.Range("H43").Value = 223842
$ws.Range("J43").Value = 223842
$ws.Range("L43").Value = 223842
$ws.Range("N43").Value = -224204

$ws.Range("H48").Value = 249684
$ws.Range("J48").Value = 249684
$ws.Range("L48").Value = 249684
$ws.Range("N48").Value = -250514

$ws.Range("H99").Value = 1527.5
$ws.Range("I99").Value = 1370
$ws.Range("K99").Value = 1370
$ws.Range("M99").Value = 128

$ws.Range("H105").Value = 1330
$ws.Range("I105").Value = 989.5
$ws.Range("J105").Value = 2011
$ws.Range("K105").Value = 989.5
$ws.Range("L105").Value = 2011
$ws.Range("M105").Value = 757.5
$ws.Range("N105").Value = -5505

$ws.Range("H134").Value = 71458.336
$ws.Range("I134").Value = 5133.9287
$ws.Range("K134").Value = 15401.7861
$ws.Range("M134").Value = -12866.7861

$ws = $wb.Worksheets("CRP")
$ws.Range("H62").Value = 4668.8
$ws.Range("I62").Value = 2481.3333
$ws.Range("K62").Value = 2481.3333
$ws.Range("M62").Value = -1857.3333

$ws.Range("H65").Value = 4668.8
$ws.Range("I65").Value = 2481.3333
$ws.Range("K65").Value = 12406.6665
$ws.Range("M65").Value = -9286.666499999999

$ws.Range("H100").Value = 65999
$ws.Range("J100").Value = 65999
$ws.Range("L100").Value = 65999
$ws.Range("N100").Value = -68163

$ws.Range("H132").Value = 1117.9166
$ws.Range("I132").Value = 889
$ws.Range("K132").Value = 2667
$ws.Range("M132").Value = -137

$ws.Range("H134").Value = 240351.47
$ws.Range("I134").Value = 2311.2683
$ws.Range("J134").Value = 10000000
$ws.Range("K134").Value = 6933.804900000001
$ws.Range("L134").Value = 30000000
$ws.Range("M134").Value = -4398.804900000001
$ws.Range("N134").Value = -30005070

$ws = $wb.Worksheets("CUL")
$ws.Range("H44").Value = 600
$ws.Range("I44").Value = 600
$ws.Range("K44").Value = 1800
$ws.Range("M44").Value = -1402

$ws.Range("H132").Value = 789633.8
$ws.Range("J132").Value = 1435286.1
$ws.Range("L132").Value = 12917574.9
$ws.Range("N132").Value = -12922634.9

$ws.Range("H137").Value = 2867.4736
$ws.Range("J137").Value = 7699.5
$ws.Range("L137").Value = 23098.5
$ws.Range("N137").Value = -33298.5

$ws = $wb.Worksheets("GSM")
$ws.Range("H102").Value = 1674.1364
$ws.Range("I102").Value = 695.94116
$ws.Range("K102").Value = 695.94116
$ws.Range("M102").Value = 926.05884

$ws.Range("H122").Value = 4798.75
$ws.Range("J122").Value = 6120
$ws.Range("L122").Value = 18360
$ws.Range("N122").Value = -23260

$ws = $wb.Worksheets("LTW")
$ws.Range("H22").Value = 1051.1936
$ws.Range("I22").Value = 503.5909
$ws.Range("J22").Value = 2389.7778
$ws.Range("K22").Value = 503.5909
$ws.Range("L22").Value = 2389.7778
$ws.Range("M22").Value = -208.5909
$ws.Range("N22").Value = -2979.7778

$ws.Range("H27").Value = 1051.1936
$ws.Range("I27").Value = 503.5909
$ws.Range("J27").Value = 2389.7778
$ws.Range("K27").Value = 503.5909
$ws.Range("L27").Value = 2389.7778
$ws.Range("M27").Value = -396.5909
$ws.Range("N27").Value = -2603.7778

$ws.Range("H127").Value = 49991
$ws.Range("J127").Value = 49991
$ws.Range("L127").Value = 49991
$ws.Range("N127").Value = -59911

$ws.Range("H132").Value = 2967.1667
$ws.Range("I132").Value = 2949.5
$ws.Range("J132").Value = 3002.5
$ws.Range("K132").Value = 8848.5
$ws.Range("L132").Value = 9007.5
$ws.Range("M132").Value = -6318.5
$ws.Range("N132").Value = -14067.5

$ws.Range("H136").Value = 1256473.8
$ws.Range("I136").Value = 1672947.8
$ws.Range("J136").Value = 7052
$ws.Range("K136").Value = 5018843.4
$ws.Range("L136").Value = 21156
$ws.Range("M136").Value = -5016293.4
$ws.Range("N136").Value = -26256

$ws = $wb.Worksheets("WVR")
$ws.Range("H62").Value = 73100.07000000001
$ws.Range("I62").Value = 206299.6
$ws.Range("K62").Value = 206299.6
$ws.Range("M62").Value = -205675.6

$ws.Range("H65").Value = 73100.07000000001
$ws.Range("I65").Value = 206299.6
$ws.Range("K65").Value = 1031498
$ws.Range("M65").Value = -1028378

$ws.Range("H100").Value = 788.375
$ws.Range("I100").Value = 788.375
$ws.Range("K100").Value = 1576.75
$ws.Range("M100").Value = -1035.75

$ws.Range("H101").Value = 29300
$ws.Range("J101").Value = 29300
$ws.Range("L101").Value = 29300
$ws.Range("N101").Value = -35790

$ws.Range("H122").Value = 47622784
$ws.Range("I122").Value = 66669850
$ws.Range("J122").Value = 5117.5
$ws.Range("K122").Value = 200009550
$ws.Range("L122").Value = 15352.5
$ws.Range("M122").Value = -200007100
$ws.Range("N122").Value = -20252.5
